$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.157.78'
$ws.Range("E2").Value = '  -2.52%  '

# Row 3
$ws.Range("D3").Value = '1.848.44'
$ws.Range("E3").Value = '  -1.43%  '

# Row 4
$ws.Range("D4").Value = "'0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").Value = "'0.6940"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.30%  '

# Row 6
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = "'238.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.89%  '

# Row 7
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").Value = "'0.3052"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.05%  '

# Row 9
$ws.Range("D9").Value = "'0.07634"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.65%  '

# Row 10
$ws.Range("D10").Value = "'23.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.03%  '

# Row 11
$ws.Range("D11").Value = "'0.08110"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.77%  '

# Row 12
$ws.Range("D12").Value = '1.853.76'
$ws.Range("E12").Value = '  -1.84%  '

# Row 13
$ws.Range("D13").Value = "'0.7238"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.51%  '

# Row 14
$ws.Range("D14").Value = "'5.179"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.94%  '

# Row 15
$ws.Range("D15").Value = "'89.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.39%  '

# Row 16
$ws.Range("D16").Value = '29.325.33'
$ws.Range("E16").Value = '  -1.94%  '

# Row 17
$ws.Range("D17").Value = "'5.788"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.36%  '

# Row 18
$ws.Range("D18").Value = "'241.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.97%  '

# Row 19
$ws.Range("D19").Value = "'0.000007718"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.58%  '

# Row 20
$ws.Range("D20").Value = "'13.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.59%  '

# Row 21
$ws.Range("E21").Value = '  +0.01%  '

# Row 22
$ws.Range("D22").Value = '2.129.34'
$ws.Range("E22").Value = '  -0.54%  '

# Row 23
$ws.Range("D23").Value = "'0.9994"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.03%  '

# Row 24
$ws.Range("D24").Value = "'7.639"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.51%  '

# Row 25
$ws.Range("D25").Value = "'9.011"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.06%  '

# Row 26
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = "'0.1461"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.74%  '

# Row 27
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = "'161.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.54%  '

# Row 28
$ws.Range("D28").Value = "'18.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.10%  '

# Row 29
$ws.Range("D29").Value = "'1.936"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.20%  '

# Row 30
$ws.Range("D30").Value = "'1.389"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.43%  '

# Row 31
$ws.Range("D31").Value = "'4.425"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.96%  '

# Row 32
$ws.Range("D32").Value = "'1.493"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.84%  '

# Row 33
$ws.Range("D33").Value = "'4.057"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.14%  '

# Row 34
$ws.Range("D34").Value = "'0.05234"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.09%  '

# Row 35
$ws.Range("D35").Value = "'1.189"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.68%  '

# Row 36
$ws.Range("D36").Value = "'0.7116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.08%  '

# Row 37
$ws.Range("D37").Value = "'1.001"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.04%  '

# Row 38
$ws.Range("D38").Value = "'2.660"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.49%  '

# Row 39
$ws.Range("D39").Value = "'0.01861"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.01%  '

# Row 40
$ws.Range("D40").Value = "'2.688"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.32%  '

# Row 41
$ws.Range("D41").Value = "'0.9159"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.22%  '

# Row 42
$ws.Range("D42").Value = "'5.952"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.09%  '

# Row 43
$ws.Range("D43").Value = "'0.4292"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.07%  '

# Row 44
$ws.Range("D44").Value = "'69.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.76%  '

# Row 45
$ws.Range("D45").Value = '1.041.78'
$ws.Range("E45").Value = '  -6.53%  '

# Row 46
$ws.Range("E46").Value = '  -0.04%  '

# Row 47
$ws.Range("D47").Value = "'102.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.54%  '

# Row 48
$ws.Range("D48").Value = "'7.225"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.08%  '

# Row 49
$ws.Range("D49").Value = '2.021.93'
$ws.Range("E49").Value = '  -0.76%  '

# Row 50
$ws.Range("D50").Value = "'1.748"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.28%  '

# Row 51
$ws.Range("D51").Value = "'9.249"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.76%  '
